$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F5").Value = 83
$ws.Range("H5").Value = 83
$ws.Range("E6").Value = 41
$ws.Range("F6").Value = 30
$ws.Range("H6").Value = 30
$ws.Range("E10").Value = 442
$ws.Range("F10").Value = 217
$ws.Range("H10").Value = 217
$ws.Range("E11").Value = 301
$ws.Range("F11").Value = 168
$ws.Range("H11").Value = 168
$ws.Range("F12").Value = 233
$ws.Range("H12").Value = 233
$ws.Range("F13").Value = 59
$ws.Range("H13").Value = 59
$ws.Range("E14").Value = 111
$ws.Range("F14").Value = 56
$ws.Range("H14").Value = 56
$ws.Range("F15").Value = 57
$ws.Range("H15").Value = 57
$ws.Range("E17").Value = 87
$ws.Range("F17").Value = 43
$ws.Range("H17").Value = 43
$ws.Range("F18").Value = 24
$ws.Range("H18").Value = 24
$ws.Range("E21").Value = 131
$ws.Range("F21").Value = 70
$ws.Range("H21").Value = 70
$ws.Range("F22").Value = 77
$ws.Range("H22").Value = 77
$ws.Range("E23").Value = 182
$ws.Range("F23").Value = 80
$ws.Range("H23").Value = 80
$ws.Range("E24").Value = 186
$ws.Range("F24").Value = 97
$ws.Range("H24").Value = 97
$ws.Range("F25").Value = 110
$ws.Range("H25").Value = 110
$ws.Range("E26").Value = 136
$ws.Range("F26").Value = 80
$ws.Range("H26").Value = 80
$ws.Range("E27").Value = 299
$ws.Range("F27").Value = 140
$ws.Range("H27").Value = 140
$ws.Range("F28").Value = 63
$ws.Range("H28").Value = 63
$ws.Range("F29").Value = 82
$ws.Range("H29").Value = 82
$ws.Range("E30").Value = 190
$ws.Range("E31").Value = 69
$ws.Range("E32").Value = 169
$ws.Range("E33").Value = 260
$ws.Range("F33").Value = 133
$ws.Range("H33").Value = 133
$ws.Range("E34").Value = 197
$ws.Range("F34").Value = 122
$ws.Range("H34").Value = 122
$ws.Range("F35").Value = 83
$ws.Range("H35").Value = 83
$ws.Range("F37").Value = 69
$ws.Range("H37").Value = 69
$ws.Range("F38").Value = 53
$ws.Range("H38").Value = 53
$ws.Range("F39").Value = 79
$ws.Range("H39").Value = 79
$ws.Range("F40").Value = 107
$ws.Range("H40").Value = 107
$ws.Range("E41").Value = 357
$ws.Range("F41").Value = 167
$ws.Range("H41").Value = 167
$ws.Range("E42").Value = 327
$ws.Range("F42").Value = 180
$ws.Range("H42").Value = 180
$ws.Range("F43").Value = 60
$ws.Range("H43").Value = 60
$ws.Range("E44").Value = 282
$ws.Range("F44").Value = 137
$ws.Range("H44").Value = 137
$ws.Range("E45").Value = 126
$ws.Range("F45").Value = 64
$ws.Range("H45").Value = 64
$ws.Range("E46").Value = 282
$ws.Range("F46").Value = 153
$ws.Range("H46").Value = 153
$ws.Range("E47").Value = 397
$ws.Range("F47").Value = 196
$ws.Range("H47").Value = 196
$ws.Range("F48").Value = 76
$ws.Range("H48").Value = 76
$ws.Range("F49").Value = 110
$ws.Range("H49").Value = 110
$ws.Range("E50").Value = 227
$ws.Range("F50").Value = 102
$ws.Range("H50").Value = 102
$ws.Range("E51").Value = 213
$ws.Range("F51").Value = 87
$ws.Range("H51").Value = 87
$ws.Range("E52").Value = 25
$ws.Range("F52").Value = 11
$ws.Range("H52").Value = 11
